$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("KPI")

# Row 10 - Meta 1 (was 0.5)
$ws.Range("C10").Value = "Meta1"
$ws.Range("J10").Value = "Meta1"
$ws.Range("Q10").Value = "Meta1"

# Row 11 - Meta 2 (was 0.8)
$ws.Range("C11").Value = "Meta2"
$ws.Range("J11").Value = "Meta2"
$ws.Range("Q11").Value = "Meta2"

# Row 12 - Meta 3 (was 1)
$ws.Range("C12").Value = "Meta3"
$ws.Range("J12").Value = "Meta3"
$ws.Range("Q12").Value = "Meta3"

# Row 13 - Meta 4 / Maeta4 (was 1.2)
$ws.Range("C13").Value = "Maeta4"
$ws.Range("J13").Value = "Meta4"
$ws.Range("Q13").Value = "Meta4"

# Row 25 - Metas de ROI headers for each of the three blocks
$ws.Range("D25").Value = "Meta1"
$ws.Range("E25").Value = "Meta2"
$ws.Range("F25").Value = "Meta3"
$ws.Range("G25").Value = "Meta4"

$ws.Range("K25").Value = "Meta1"
$ws.Range("L25").Value = "Meta2"
$ws.Range("M25").Value = "Meta3"
$ws.Range("N25").Value = "Meta4"

$ws.Range("R25").Value = "Meta1"
$ws.Range("S25").Value = "Meta2"
$ws.Range("T25").Value = "Meta3"
$ws.Range("U25").Value = "Meta4"

# Move the active selection on the KPI sheet (cosmetic, matches the new layout)
$ws.Activate()
$ws.Range("B46").Select()
